$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "Juan Luis" SmartScore cells were saved as text by the Streamlit
#     app; re-entering the same figures through Excel stores them as real
#     numbers (matches the upstream fix committed alongside the new row). ---
$ws.Range("G3").Value = 0.572
$ws.Range("J3").Value = 0.514
$ws.Range("M3").Value = 0.409
$ws.Range("P3").Value = 0.845
$ws.Range("S3").Value = 0.618
$ws.Range("V3").Value = 0.602
$ws.Range("Y3").Value = 0.769
$ws.Range("AB3").Value = 0.503
$ws.Range("AE3").Value = 0.423

# --- Row 4: new SmartScore submission from Rosa Linda ---
$ws.Range("A4").Value = 'Rosa Linda'
$ws.Range("B4").Value = 24
$ws.Range("C4").Value = 'Femenino'
$ws.Range("D4").Value = '2025-10-28 05:43:11'
$ws.Range("E4").Value = '{
  "portion": 0.4,
  "diet": 0.7142857142857143,
  "salt": 0.4,
  "fat": 1.0,
  "natural": 1.0,
  "convenience": 0.8,
  "price": 0.8
}'
$ws.Range("F4").Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '0.572'
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("I4").Value = 'Maruchan Ramen Sabor Pollo'
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = '0.488'
$ws.Range("J4").ClearFormats()
$ws.Range("K4").Value = 'Sabor clásico, económico, alto en sodio, no saludable, nostálgico'
$ws.Range("L4").Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = '0.412'
$ws.Range("M4").ClearFormats()
$ws.Range("N4").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("O4").Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = '0.636'
$ws.Range("P4").ClearFormats()
$ws.Range("Q4").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("R4").Value = 'Annie’s Shells & White Cheddar'
$ws.Range("S4").NumberFormat = "@"
$ws.Range("S4").Value = '0.578'
$ws.Range("S4").ClearFormats()
$ws.Range("T4").Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("U4").Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range("V4").NumberFormat = "@"
$ws.Range("V4").Value = '0.541'
$ws.Range("V4").ClearFormats()
$ws.Range("W4").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("X4").Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = '0.738'
$ws.Range("Y4").ClearFormats()
$ws.Range("Z4").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AA4").Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = '0.614'
$ws.Range("AB4").ClearFormats()
$ws.Range("AC4").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
$ws.Range("AD4").Value = 'Jack Link’s Beef Jerky Original'
$ws.Range("AE4").NumberFormat = "@"
$ws.Range("AE4").Value = '0.599'
$ws.Range("AE4").ClearFormats()
$ws.Range("AF4").Value = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'

# Undo the automatic row-height bump Excel applies for the multi-line
# JSON "Pesos" cell (E4) so the new row keeps the default row height.
$ws.Rows.Item(4).AutoFit()
